$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 63.91118233333333
$ws.Range("H2").Value = 191.733547
$ws.Range("I2").Value = 0.4067926910433548
$ws.Range("J2").Value = 0.4067926910433549
$ws.Range("M2").Value = 0.2374196666666667
$ws.Range("N2").Value = 0.7122590000000001
$ws.Range("O2").Value = 0.04883833429957388
$ws.Range("P2").Value = 0.04883833429957388
$ws.Range("Q2").Value = 15.17377160585256
$ws.Range("R2").Value = 136.563944452673
$ws.Range("S2").Value = 0.01986707743579863
$ws.Range("T2").Value = 0.01986707743579864

$ws.Range("G3").Value = 63.91118233333333
$ws.Range("H3").Value = 191.733547
$ws.Range("I3").Value = 0.4067926910433548
$ws.Range("J3").Value = 0.4067926910433549
$ws.Range("O3").Value = 0.1945563001683692
$ws.Range("P3").Value = 0.1945563001683693
$ws.Range("Q3").Value = 60.447451895596
$ws.Range("R3").Value = 544.027067060364
$ws.Range("S3").Value = 0.07914408090492962
$ws.Range("T3").Value = 0.07914408090492965

$ws.Range("G4").Value = 63.91118233333333
$ws.Range("H4").Value = 191.733547
$ws.Range("I4").Value = 0.4067926910433548
$ws.Range("J4").Value = 0.4067926910433549
$ws.Range("M4").Value = 3.678114666666667
$ws.Range("N4").Value = 11.034344
$ws.Range("O4").Value = 0.7566053655320568
$ws.Range("P4").Value = 0.7566053655320568
$ws.Range("Q4").Value = 235.0726571042409
$ws.Range("R4").Value = 2115.653913938168
$ws.Range("S4").Value = 0.3077815327026265
$ws.Range("T4").Value = 0.3077815327026265

$ws.Range("I5").Value = 0.3656254573230189
$ws.Range("J5").Value = 0.365625457323019
$ws.Range("M5").Value = 0.2374196666666667
$ws.Range("N5").Value = 0.7122590000000001
$ws.Range("O5").Value = 0.04883833429957388
$ws.Range("P5").Value = 0.04883833429957388
$ws.Range("Q5").Value = 13.6381928802
$ws.Range("R5").Value = 122.7437359218
$ws.Range("S5").Value = 0.01785653831317618
$ws.Range("T5").Value = 0.01785653831317618

$ws.Range("I6").Value = 0.3656254573230189
$ws.Range("J6").Value = 0.365625457323019
$ws.Range("O6").Value = 0.1945563001683692
$ws.Range("P6").Value = 0.1945563001683693
$ws.Range("R6").Value = 488.9717774424
$ws.Range("S6").Value = 0.07113473622413455
$ws.Range("T6").Value = 0.07113473622413458

$ws.Range("I7").Value = 0.3656254573230189
$ws.Range("J7").Value = 0.365625457323019
$ws.Range("M7").Value = 3.678114666666667
$ws.Range("N7").Value = 11.034344
$ws.Range("O7").Value = 0.7566053655320568
$ws.Range("P7").Value = 0.7566053655320568
$ws.Range("Q7").Value = 211.2834120432
$ws.Range("R7").Value = 1901.5507083888
$ws.Range("S7").Value = 0.2766341827857082
$ws.Range("T7").Value = 0.2766341827857082

$ws.Range("G8").Value = 35.755375
$ws.Range("H8").Value = 107.266125
$ws.Range("I8").Value = 0.2275818516336261
$ws.Range("J8").Value = 0.2275818516336262
$ws.Range("M8").Value = 0.2374196666666667
$ws.Range("N8").Value = 0.7122590000000001
$ws.Range("O8").Value = 0.04883833429957388
$ws.Range("P8").Value = 0.04883833429957388
$ws.Range("Q8").Value = 8.489029214041668
$ws.Range("R8").Value = 76.40126292637501
$ws.Range("S8").Value = 0.01111471855059906
$ws.Range("T8").Value = 0.01111471855059906

$ws.Range("G9").Value = 35.755375
$ws.Range("H9").Value = 107.266125
$ws.Range("I9").Value = 0.2275818516336261
$ws.Range("J9").Value = 0.2275818516336262
$ws.Range("O9").Value = 0.1945563001683692
$ws.Range("P9").Value = 0.1945563001683693
$ws.Range("Q9").Value = 33.8175766965
$ws.Range("R9").Value = 304.3581902685
$ws.Range("S9").Value = 0.04427748303930503
$ws.Range("T9").Value = 0.04427748303930505

$ws.Range("G10").Value = 35.755375
$ws.Range("H10").Value = 107.266125
$ws.Range("I10").Value = 0.2275818516336261
$ws.Range("J10").Value = 0.2275818516336262
$ws.Range("M10").Value = 3.678114666666667
$ws.Range("N10").Value = 11.034344
$ws.Range("O10").Value = 0.7566053655320568
$ws.Range("P10").Value = 0.7566053655320568
$ws.Range("Q10").Value = 131.5123691996667
$ws.Range("R10").Value = 1183.611322797
$ws.Range("S10").Value = 0.172189650043722
$ws.Range("T10").Value = 0.172189650043722

Write-Output "done"
